$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "maa://24702 (94.1), maa://25390 (97.41), maa://36681 (91.94)"
$ws.Range("AE2").Value = "***maa://21730 (17.19), maa://25251 (92.0), ***maa://39501 (20.0), *maa://36675 (60.0)"
$ws.Range("C3").Value = "maa://36987 (95.65), maa://40192 (100.0), maa://39849 (100.0)"
$ws.Range("G3").Value = "maa://21247 (98.26), *maa://22748 (75.0)"
$ws.Range("K3").Value = "*maa://22880 (70.2), maa://20276 (82.35), *maa://22749 (62.5)"
$ws.Range("O3").Value = "maa://21249 (95.5), maa://26254 (95.24)"
$ws.Range("C4").Value = "maa://24632 (93.18), **maa://24303 (36.36), maa://22499 (85.71), maa://22746 (100.0)"
$ws.Range("S4").Value = "maa://32509 (98.72), maa://22754 (91.67), maa://27295 (80.39), *maa://21746 (55.81), *maa://31008 (78.05)"
$ws.Range("W4").Value = "**maa://32495 (47.9), ***maa://31785 (15.74), ***maa://36683 (26.67)"
$ws.Range("C5").Value = "maa://21245 (82.54), maa://22744 (82.61)"
$ws.Range("K7").Value = "maa://28624 (91.3), maa://24957 (97.3)"
$ws.Range("W7").Value = "maa://22399 (94.62), *maa://22758 (71.43)"
$ws.Range("AE7").Value = "*maa://26191 (70.42), *maa://36671 (73.17)"
$ws.Range("W8").Value = "maa://21411 (96.25)"
$ws.Range("AA9").Value = "maa://28711 (87.65), ***maa://22740 (5.88), **maa://27377 (46.15), ***maa://25174 (20.0), **maa://39938 (50.0), maa://40166 (100.0)"
$ws.Range("AE9").Value = "maa://26206 (90.0), **maa://22865 (45.65)"
$ws.Range("C10").Value = "***maa://25695 (19.41), **maa://32237 (38.89), ***maa://34206 (14.29), ***maa://39951 (15.79), ***maa://39243 (25.0)"
$ws.Range("S10").Value = "maa://27395 (97.24), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("W10").Value = "maa://22301 (97.33), maa://22726 (100.0)"
$ws.Range("AE10").Value = "*maa://25021 (56.52), *maa://22733 (58.62), maa://22761 (100.0)"
$ws.Range("C11").Value = "maa://36707 (99.63)"
$ws.Range("S11").Value = "maa://22747 (94.96), maa://22501 (98.08)"
$ws.Range("W11").Value = "maa://36713 (97.78)"
$ws.Range("G12").Value = "maa://21867 (89.86)"
$ws.Range("W12").Value = "maa://22753 (91.67), *maa://21485 (76.56), maa://37962 (81.25)"
$ws.Range("AA12").Value = "maa://23669 (95.82), maa://36677 (94.74), maa://39872 (81.82)"
$ws.Range("C13").Value = "maa://24999 (91.36), maa://36673 (91.8), maa://25001 (85.51)"
$ws.Range("G13").Value = "*maa://21248 (75.5), **maa://22728 (47.62)"
$ws.Range("O13").Value = "maa://22676 (91.75), *maa://22583 (75.86), *maa://22500 (55.81)"
$ws.Range("W13").Value = "*maa://34957 (79.49), *maa://22768 (53.33)"
$ws.Range("AE13").Value = "**maa://22737 (30.83), maa://39883 (90.48), *maa://39885 (73.68)"
$ws.Range("K14").Value = "maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (92.11)"
$ws.Range("C15").Value = "*maa://22743 (76.05), maa://22734 (83.33), *maa://30808 (63.64), ***maa://36048 (13.33)"
$ws.Range("G15").Value = "maa://24304 (88.27), maa://21478 (91.18)"
$ws.Range("O15").Value = "maa://24762 (89.51), *maa://22727 (70.0)"
$ws.Range("AE15").Value = "maa://21364 (80.55), *maa://22766 (73.0), *maa://36666 (76.67)"
$ws.Range("C16").Value = "maa://21441 (96.15), maa://36679 (93.33), maa://37650 (95.0)"
$ws.Range("W16").Value = "maa://28501 (97.37), maa://28051 (95.83)"
$ws.Range("AA16").Value = "maa://26228 (96.0)"
$ws.Range("AE16").Value = "*maa://23911 (61.54), maa://27755 (91.55)"
$ws.Range("C18").Value = "maa://24570 (96.43)"
$ws.Range("G18").Value = "maa://24421 (90.34)"
$ws.Range("W18").Value = "maa://21917 (97.4), maa://22741 (83.33)"
$ws.Range("O20").Value = "maa://37442 (96.3)"
$ws.Range("K21").Value = "maa://31731 (95.0)"
$ws.Range("AA21").Value = "*maa://21443 (78.64), ***maa://23820 (29.63)"
$ws.Range("AE21").Value = "maa://22524 (94.22), *maa://22432 (75.47)"
$ws.Range("K22").Value = "maa://27127 (87.36), *maa://22751 (77.05)"
$ws.Range("W22").Value = "maa://21282 (98.8), *maa://37649 (64.71)"
$ws.Range("K23").Value = "maa://39756 (92.03), maa://39875 (94.74)"
$ws.Range("O23").Value = "maa://30587 (91.62), *maa://29748 (74.8), ***maa://29785 (15.15), *maa://37566 (76.47)"
$ws.Range("W24").Value = "maa://23504 (93.02), maa://29988 (85.86), **maa://22892 (40.43), *maa://25141 (77.31), *maa://36663 (79.25), ***maa://22815 (23.08)"
$ws.Range("AE24").Value = "maa://22523 (84.86), *maa://36672 (75.61), maa://29910 (95.74), **maa://21440 (34.55)"
$ws.Range("G25").Value = "*maa://29063 (77.95), *maa://25311 (74.73), ***maa://22725 (4.84)"
$ws.Range("K25").Value = "maa://24378 (88.24)"
$ws.Range("AA25").Value = "maa://31215 (85.9), *maa://24516 (80.0), maa://26001 (88.89)"
$ws.Range("G26").Value = "maa://24913 (90.91)"
$ws.Range("AE26").Value = "maa://30511 (83.87), *maa://29760 (54.55)"
$ws.Range("C28").Value = "maa://24465 (90.3), maa://25725 (81.82)"
$ws.Range("W28").Value = "maa://39929 (85.21), ***maa://39723 (15.15)"
$ws.Range("AE28").Value = "maa://36660 (94.02), *maa://36701 (64.0)"
$ws.Range("G29").Value = "*maa://25175 (70.45)"
$ws.Range("K29").Value = "maa://28432 (93.33), *maa://28440 (72.5), maa://31400 (100.0), *maa://28650 (66.67)"
$ws.Range("O29").Value = "*maa://23168 (55.77), **maa://30050 (42.11)"
$ws.Range("AE29").Value = "*maa://24080 (68.25), ***maa://34960 (9.09)"
$ws.Range("K30").Value = "maa://30442 (94.12)"
$ws.Range("K31").Value = "maa://35926 (93.42), maa://36258 (80.26)"
$ws.Range("S32").Value = "maa://41108 (88.89), maa://41238 (93.33)"
$ws.Range("O33").Value = "*maa://21956 (78.91), maa://22730 (82.14)"
$ws.Range("K35").Value = "maa://41296 (96.15)"
$ws.Range("AE38").Value = "maa://36697 (86.78)"
$ws.Range("G39").Value = "maa://25199 (86.11), maa://36670 (88.33), maa://30434 (88.89), ***maa://25036 (16.0)"
$ws.Range("O40").Value = "maa://23278 (96.19), maa://21386 (95.63), maa://36664 (92.11)"
$ws.Range("G43").Value = "maa://22525 (92.56), maa://21284 (82.93)"
$ws.Range("G45").Value = "maa://21229 (85.47), maa://30807 (94.92), *maa://22767 (52.94), ***maa://20796 (13.79)"
$ws.Range("G46").Value = "maa://35931 (92.34)"
$ws.Range("G47").Value = "maa://27410 (95.79), maa://29661 (97.6), maa://28038 (84.62)"
$ws.Range("G53").Value = "maa://32534 (93.16), **maa://32434 (36.36)"
$ws.Range("G55").Value = "maa://32532 (92.13)"
$ws.Range("G57").Value = "maa://25176 (97.73)"
$ws.Range("G58").Value = "*maa://37964 (58.82)"
